$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B for "Week_Start_Date" (shifts ASIN..is_holiday_week from B:I to C:J)
$ws.Columns.Item(2).Insert()

# Force column B to be stored as text so the week-start dates are not auto-converted to date serials
$ws.Columns.Item(2).NumberFormat = "@"

# Header row
$ws.Range("B1").Value = "Week_Start_Date"

# Data rows: update Week label, set Week_Start_Date, and refresh MyForecast (column D after insert)
$ws.Range("A2").Value = "W1"
$ws.Range("B2").Value = "2025-01-05"
$ws.Range("D2").Value = 100
$ws.Range("A3").Value = "W2"
$ws.Range("B3").Value = "2025-01-12"
$ws.Range("D3").Value = 95
$ws.Range("A4").Value = "W3"
$ws.Range("B4").Value = "2025-01-19"
$ws.Range("D4").Value = 101
$ws.Range("A5").Value = "W4"
$ws.Range("B5").Value = "2025-01-26"
$ws.Range("D5").Value = 99
$ws.Range("A6").Value = "W5"
$ws.Range("B6").Value = "2025-02-02"
$ws.Range("D6").Value = 96
$ws.Range("A7").Value = "W6"
$ws.Range("B7").Value = "2025-02-09"
$ws.Range("D7").Value = 93
$ws.Range("A8").Value = "W7"
$ws.Range("B8").Value = "2025-02-16"
$ws.Range("D8").Value = 98
$ws.Range("A9").Value = "W8"
$ws.Range("B9").Value = "2025-02-23"
$ws.Range("D9").Value = 97
$ws.Range("A10").Value = "W9"
$ws.Range("B10").Value = "2025-03-02"
$ws.Range("D10").Value = 92
$ws.Range("A11").Value = "W10"
$ws.Range("B11").Value = "2025-03-09"
$ws.Range("D11").Value = 88
$ws.Range("A12").Value = "W11"
$ws.Range("B12").Value = "2025-03-16"
$ws.Range("D12").Value = 86
$ws.Range("A13").Value = "W12"
$ws.Range("B13").Value = "2025-03-23"
$ws.Range("D13").Value = 82
$ws.Range("A14").Value = "W13"
$ws.Range("B14").Value = "2025-03-30"
$ws.Range("D14").Value = 81
$ws.Range("A15").Value = "W14"
$ws.Range("B15").Value = "2025-04-06"
$ws.Range("D15").Value = 85
$ws.Range("A16").Value = "W15"
$ws.Range("B16").Value = "2025-04-13"
$ws.Range("D16").Value = 92
$ws.Range("A17").Value = "W16"
$ws.Range("B17").Value = "2025-04-20"
$ws.Range("D17").Value = 89

# is_holiday_week column (J) becomes boolean FALSE instead of numeric 0
$ws.Range("J2:J17").Value = $false

$ws.Range("A1").Select()
